# Generate Report for Archive
#
# The localization status changed from "Ready for handoff" to
# "In Translation" for the three sample rows, on all three sheets
# (Overview, zh-cn, de-de). Updating the text makes the "Status"-ish
# columns narrower, so their column widths are refreshed to fit the
# new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"

# --- zh-cn sheet: column C (Status) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C4").Value = "In Translation"

# --- de-de sheet: column C (Status) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C4").Value = "In Translation"

# Resize the affected columns to fit the new, shorter status text.
$wsOverview.Range("E:E").ColumnWidth = 12.5
$wsOverview.Range("F:F").ColumnWidth = 12.5
$wsZhCn.Range("C:C").ColumnWidth = 12.5
$wsDeDe.Range("C:C").ColumnWidth = 12.5
